# Auto-generated edit script applying market-price/profit refresh diffs
# to the Halicarnassus_Profits workbook, matching the supplied OOXML diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3849.3
$ws.Range("I116").Value = 3882.1667
$ws.Range("K116").Value = 3882.1667
$ws.Range("M116").Value = -440.1667000000002
$ws.Range("H127").Value = 1633.3334
$ws.Range("I127").Value = 1633.3334
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 4900.0002
$ws.Range("L127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = 59.9997999999996

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 10000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 10000
$ws.Range("K27").Value = 0
$ws.Range("L27").ClearContents()
$ws.Range("M27").Value = 10000
$ws.Range("N27").Value = -10368
$ws.Range("H37").Value = 19225
$ws.Range("J37").Value = 18450
$ws.Range("L37").Value = 18450
$ws.Range("N37").Value = -18996
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("N40").Value = 0
$ws.Range("H45").Value = 2439.0625
$ws.Range("I45").Value = 1140
$ws.Range("K45").Value = 1140
$ws.Range("M45").Value = -763
$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 2000
$ws.Range("K61").Value = 2000
$ws.Range("M61").Value = -1788
$ws.Range("H63").Value = 1878.1666
$ws.Range("I63").Value = 1878.1666
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1878.1666
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -1192.1666
$ws.Range("H66").Value = 1878.1666
$ws.Range("I66").Value = 1878.1666
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 9390.833000000001
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -5958.833000000001
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 717.36365
$ws.Range("I80").Value = 324.5
$ws.Range("K80").Value = 324.5
$ws.Range("M80").Value = 673.5
$ws.Range("H82").Value = 21875
$ws.Range("I82").Value = 6672
$ws.Range("K82").Value = 6672
$ws.Range("M82").Value = -6289
$ws.Range("H83").Value = 717.36365
$ws.Range("I83").Value = 324.5
$ws.Range("K83").Value = 1622.5
$ws.Range("M83").Value = 3369.5
$ws.Range("H85").Value = 21875
$ws.Range("I85").Value = 6672
$ws.Range("K85").Value = 6672
$ws.Range("M85").Value = -5346
$ws.Range("H105").Value = 1478.7858
$ws.Range("I105").Value = 1206.2222
$ws.Range("K105").Value = 1206.2222
$ws.Range("M105").Value = 540.7778000000001
$ws.Range("H107").Value = 9169.429
$ws.Range("I107").Value = 9056.714
$ws.Range("J107").Value = 9282.143
$ws.Range("K107").Value = 9056.714
$ws.Range("L107").Value = 9282.143
$ws.Range("M107").Value = -7136.714
$ws.Range("N107").Value = -13122.143

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 2121.5
$ws.Range("I10").Value = 829.7143
$ws.Range("J10").Value = 5135.6665
$ws.Range("K10").Value = 829.7143
$ws.Range("L10").Value = 5135.6665
$ws.Range("M10").Value = -690.7143
$ws.Range("N10").Value = -5413.6665
$ws.Range("H19").Value = 15333448
$ws.Range("I19").Value = 23000072
$ws.Range("K19").Value = 23000072
$ws.Range("M19").Value = -22999902
$ws.Range("H24").Value = 15333448
$ws.Range("I24").Value = 23000072
$ws.Range("K24").Value = 23000072
$ws.Range("M24").Value = -22999902
$ws.Range("H31").Value = 6168.5293
$ws.Range("I31").Value = 1858.125
$ws.Range("K31").Value = 1858.125
$ws.Range("M31").Value = -1563.125
$ws.Range("H34").Value = 6168.5293
$ws.Range("I34").Value = 1858.125
$ws.Range("K34").Value = 1858.125
$ws.Range("M34").Value = -1656.125
$ws.Range("H35").Value = 1562.1666
$ws.Range("I35").Value = 1562.1666
$ws.Range("K35").Value = 1562.1666
$ws.Range("M35").Value = -1268.1666
$ws.Range("H39").Value = 7994
$ws.Range("I39").Value = 2105
$ws.Range("J39").Value = 9957
$ws.Range("K39").Value = 2105
$ws.Range("L39").Value = 9957
$ws.Range("M39").Value = -1714
$ws.Range("N39").Value = -10739
$ws.Range("H49").Value = 7994
$ws.Range("I49").Value = 2105
$ws.Range("J49").Value = 9957
$ws.Range("K49").Value = 2105
$ws.Range("L49").Value = 9957
$ws.Range("M49").Value = -1923
$ws.Range("N49").Value = -10321
$ws.Range("H94").Value = 8503.25
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 11004.333
$ws.Range("K94").Value = 1000
$ws.Range("L94").Value = 11004.333
$ws.Range("M94").Value = -549
$ws.Range("N94").Value = -11906.333

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 623.2
$ws.Range("I7").Value = 338.66666
$ws.Range("K7").Value = 1015.99998
$ws.Range("M7").Value = -903.9999799999999
$ws.Range("H23").Value = 109.44444
$ws.Range("I23").Value = 31.75
$ws.Range("K23").Value = 95.25
$ws.Range("M23").Value = 139.75
$ws.Range("H114").Value = 1040.6364
$ws.Range("I114").Value = 952.4286
$ws.Range("J114").Value = 1195
$ws.Range("K114").Value = 2857.2858
$ws.Range("L114").Value = 3585
$ws.Range("M114").Value = 396.7142000000003
$ws.Range("N114").Value = -10093
$ws.Range("H117").Value = 490.0435
$ws.Range("I117").Value = 144.5
$ws.Range("J117").Value = 522.9524
$ws.Range("K117").Value = 433.5
$ws.Range("L117").Value = 1568.8572
$ws.Range("M117").Value = 3008.5
$ws.Range("N117").Value = -8452.8572
$ws.Range("H121").Value = 1533.25
$ws.Range("I121").Value = 359.8
$ws.Range("J121").Value = 2371.4285
$ws.Range("K121").Value = 1079.4
$ws.Range("L121").Value = 7114.2855
$ws.Range("M121").Value = 230.5999999999999
$ws.Range("N121").Value = -9734.2855
$ws.Range("H131").Value = 1675.3
$ws.Range("I131").Value = 719.6
$ws.Range("J131").Value = 2631
$ws.Range("K131").Value = 2158.8
$ws.Range("L131").Value = 7893
$ws.Range("M131").Value = 2881.2
$ws.Range("N131").Value = -17973

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10000
$ws.Range("I5").Value = 10000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 10000
$ws.Range("L5").ClearContents()
$ws.Range("M5").Value = -9888
$ws.Range("N5").Value = 0
$ws.Range("H35").Value = 27748.5
$ws.Range("J35").Value = 27748.5
$ws.Range("L35").Value = 27748.5
$ws.Range("N35").Value = -28344.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 412.16666
$ws.Range("I16").Value = 412.16666
$ws.Range("K16").Value = 412.16666
$ws.Range("M16").Value = -242.16666
$ws.Range("H32").Value = 12671.333
$ws.Range("I32").Value = 12671.333
$ws.Range("K32").Value = 12671.333
$ws.Range("M32").Value = -12354.333
$ws.Range("H55").Value = 1048.25
$ws.Range("I55").Value = 1003.55554
$ws.Range("K55").Value = 1003.55554
$ws.Range("M55").Value = -830.55554
$ws.Range("H61").Value = 5206.8184
$ws.Range("I61").Value = 2999.75
$ws.Range("J61").Value = 6468
$ws.Range("K61").Value = 2999.75
$ws.Range("L61").Value = 6468
$ws.Range("M61").Value = -2797.75
$ws.Range("N61").Value = -6872
$ws.Range("H113").Value = 5206.8184
$ws.Range("I113").Value = 2999.75
$ws.Range("J113").Value = 6468
$ws.Range("K113").Value = 2999.75
$ws.Range("L113").Value = 6468
$ws.Range("M113").Value = -829.75
$ws.Range("N113").Value = -10808
